# Append: 2026-02-08 02:51 JST
# Update the "取得日時" (retrieved datetime) timestamps for the newly
# appended rows (2-6) on the "ランサーズ" sheet from 02:23:10 to 02:51:03.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-08 02:51:03"

$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
